$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F-column "想去人数" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3150
$ws1.Range("F3").Value = 542
$ws1.Range("F4").Value = 1107
$ws1.Range("F5").Value = 95
$ws1.Range("F6").Value = 49
$ws1.Range("F7").Value = 273
$ws1.Range("F9").Value = 1135
$ws1.Range("F10").Value = 15927
$ws1.Range("F11").Value = 254
$ws1.Range("F12").Value = 186
$ws1.Range("F13").Value = 1029
$ws1.Range("F14").Value = 6236
$ws1.Range("F15").Value = 628
$ws1.Range("F16").Value = 116
$ws1.Range("F17").Value = 71
$ws1.Range("F18").Value = 12
$ws1.Range("F19").Value = 121
$ws1.Range("F20").Value = 1266
$ws1.Range("F21").Value = 31
$ws1.Range("F27").Value = 878
$ws1.Range("F29").Value = 5015
$ws1.Range("F30").Value = 493
$ws1.Range("F31").Value = 11141
$ws1.Range("F32").Value = 1240
$ws1.Range("F36").Value = 3813

# Sheet "全部类型" (sheet4): same events, rows shifted by +1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3150
$ws4.Range("F4").Value = 542
$ws4.Range("F5").Value = 1107
$ws4.Range("F6").Value = 95
$ws4.Range("F7").Value = 49
$ws4.Range("F8").Value = 273
$ws4.Range("F10").Value = 1135
$ws4.Range("F11").Value = 15927
$ws4.Range("F12").Value = 254
$ws4.Range("F13").Value = 186
$ws4.Range("F14").Value = 1029
$ws4.Range("F15").Value = 6236
$ws4.Range("F16").Value = 628
$ws4.Range("F17").Value = 116
$ws4.Range("F18").Value = 71
$ws4.Range("F19").Value = 12
$ws4.Range("F20").Value = 121
$ws4.Range("F21").Value = 1266
$ws4.Range("F22").Value = 31
$ws4.Range("F28").Value = 878
$ws4.Range("F30").Value = 5015
$ws4.Range("F31").Value = 493
$ws4.Range("F32").Value = 11141
$ws4.Range("F33").Value = 1240
$ws4.Range("F37").Value = 3813
